# The "Wrong template, one camera, channel 1" sample workbook used to store
# only the bare file names of the two images ("Img1_C1.tif" / "Img2_C1.tif").
# The finder now needs the real, absolute path on disk, so the two cells
# (and the shared-string table backing them) are updated to hold the full
# path instead of just the file name.
#
# In the original edit the two now-empty rows above them (rows 4 and 5)
# were also removed, which shifts the remaining rows up by two (old row 6
# -> new row 4, old row 7 -> new row 5, old row 8 -> new row 6) and shrinks
# the used range from A1:D8 down to A1:D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$img1Path = "/home/masoud/Documents/four-polar/fourPolar-io/src/test/resources/fr/fresnel/fourPolar/io/imageSet/acquisition/sample/finders/excel/OneCamera/Img1_C1.tif"
$img2Path = "/home/masoud/Documents/four-polar/fourPolar-io/src/test/resources/fr/fresnel/fourPolar/io/imageSet/acquisition/sample/finders/excel/OneCamera/Img2_C1.tif"

# Remove the blank rows 4 and 5, shifting rows 6-8 up to rows 4-6.
$ws.Rows("4:5").Delete()

# Rows 4-6 now hold (old row 6) "Pol0_45_90_135", (old row 7) the first
# image file, and (old row 8) the second image file in both A and B.
# Replace the file-name-only text with the full path on disk.
$ws.Range("A5").Value = $img1Path
$ws.Range("A6").Value = $img2Path
$ws.Range("B6").Value = $img2Path

# Mirror the author's resulting selection of B6.
$ws.Range("B6").Select() | Out-Null
